# Prepend a new day's price row to the "商品价格数据" (commodity price data)
# sheet. This mirrors the automated daily-update commit that inserts a new
# row right under the header with the latest date, pushing all the
# historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row just below the header (row 1), shifting all
# existing data rows (old row 2 -> row 3, old row 3 -> row 4, ...) down.
$ws.Rows.Item(2).Insert()

# Excel would otherwise infer a date value/format for a "YYYY-MM-DD"
# looking string; force the cell to text first so it stays a plain string
# like every other date cell in column A.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2026-01-11"
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610

# The row-insert operation copies formatting from the row above (the bold,
# bordered header), and forcing the NumberFormat added a style too. Strip
# all of that so the new row matches the unstyled look of the other data
# rows, exactly like the source data.
$ws.Rows.Item(2).ClearFormats()
